{"js": "// Update the two-digit multiplication answer table: replace each old\n// \"A\u00d7B=C\" cell value with its new value. Every old value is unique in\n// the document, so a direct text search + in-place replace keeps the\n// existing run formatting (font/size) untouched.\nconst replacements = [\n  [\"39\u00d735=1365\", \"95\u00d714=1330\"],\n  [\"72\u00d769=4968\", \"63\u00d796=6048\"],\n  [\"81\u00d736=2916\", \"63\u00d778=4914\"],\n  [\"70\u00d760=4200\", \"49\u00d758=2842\"],\n  [\"57\u00d727=1539\", \"90\u00d778=7020\"],\n  [\"66\u00d727=1782\", \"43\u00d728=1204\"],\n  [\"65\u00d793=6045\", \"47\u00d786=4042\"],\n  [\"23\u00d766=1518\", \"80\u00d774=5920\"],\n  [\"59\u00d726=1534\", \"79\u00d742=3318\"],\n  [\"53\u00d762=3286\", \"33\u00d731=1023\"],\n  [\"50\u00d793=4650\", \"30\u00d777=2310\"],\n  [\"50\u00d736=1800\", \"65\u00d731=2015\"],\n  [\"58\u00d777=4466\", \"15\u00d716=240\"],\n  [\"25\u00d722=550\", \"13\u00d739=507\"],\n  [\"11\u00d727=297\", \"12\u00d733=396\"],\n  [\"51\u00d763=3213\", \"61\u00d753=3233\"],\n  [\"47\u00d756=2632\", \"88\u00d720=1760\"],\n  [\"31\u00d750=1550\", \"65\u00d736=2340\"],\n  [\"19\u00d765=1235\", \"75\u00d771=5325\"],\n  [\"49\u00d731=1519\", \"92\u00d725=2300\"],\n  [\"99\u00d723=2277\", \"24\u00d732=768\"],\n  [\"38\u00d788=3344\", \"37\u00d796=3552\"],\n  [\"16\u00d759=944\", \"74\u00d750=3700\"],\n  [\"42\u00d767=2814\", \"89\u00d779=7031\"],\n  [\"56\u00d712=672\", \"44\u00d717=748\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication answer table: replace each old\n# \"A\u00d7B=C\" cell value with its new value. Every old value is unique in\n# the document, so Find/Replace on the whole document content keeps\n# the existing run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"39\u00d735=1365\", \"95\u00d714=1330\"),\n    @(\"72\u00d769=4968\", \"63\u00d796=6048\"),\n    @(\"81\u00d736=2916\", \"63\u00d778=4914\"),\n    @(\"70\u00d760=4200\", \"49\u00d758=2842\"),\n    @(\"57\u00d727=1539\", \"90\u00d778=7020\"),\n    @(\"66\u00d727=1782\", \"43\u00d728=1204\"),\n    @(\"65\u00d793=6045\", \"47\u00d786=4042\"),\n    @(\"23\u00d766=1518\", \"80\u00d774=5920\"),\n    @(\"59\u00d726=1534\", \"79\u00d742=3318\"),\n    @(\"53\u00d762=3286\", \"33\u00d731=1023\"),\n    @(\"50\u00d793=4650\", \"30\u00d777=2310\"),\n    @(\"50\u00d736=1800\", \"65\u00d731=2015\"),\n    @(\"58\u00d777=4466\", \"15\u00d716=240\"),\n    @(\"25\u00d722=550\", \"13\u00d739=507\"),\n    @(\"11\u00d727=297\", \"12\u00d733=396\"),\n    @(\"51\u00d763=3213\", \"61\u00d753=3233\"),\n    @(\"47\u00d756=2632\", \"88\u00d720=1760\"),\n    @(\"31\u00d750=1550\", \"65\u00d736=2340\"),\n    @(\"19\u00d765=1235\", \"75\u00d771=5325\"),\n    @(\"49\u00d731=1519\", \"92\u00d725=2300\"),\n    @(\"99\u00d723=2277\", \"24\u00d732=768\"),\n    @(\"38\u00d788=3344\", \"37\u00d796=3552\"),\n    @(\"16\u00d759=944\", \"74\u00d750=3700\"),\n    @(\"42\u00d767=2814\", \"89\u00d779=7031\"),\n    @(\"56\u00d712=672\", \"44\u00d717=748\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $r = $d.Content\n    $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
